# Natmi following Dr Hou advice
#
# The ligand/receptor-expressing cell counts for this cluster pair table were
# recomputed upstream (ECs/FAPs/sCs now each contribute 3 expressing cells
# instead of 1), which changes the derived average/total expression values
# and every specificity score that is derived from them. Apply the new
# values for rows 2-10 (columns E, G, H, I, J, K, M, N, O, P, Q, R, S, T);
# columns A-D, F and L are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;  E = 3; G = 91.47356466666666;  H = 274.420694;          I = 0.7914295280104694; J = 0.7914295280104694; K = 3; M = 21.870458;           N = 65.61137400000001;  O = 0.9662106525726075;  P = 0.9662106525726075;  Q = 2000.568754152617;   R = 18005.11878737356;  S = 0.7646876407242263;  T = 0.7646876407242263  }
    @{ Row = 3;  E = 3; G = 91.47356466666666;  H = 274.420694;          I = 0.7914295280104694; J = 0.7914295280104694; K = 3; M = 0.1926236666666667;  N = 0.577871;           O = 0.008509882997005752; P = 0.008509882997005752; Q = 17.61997342916378;   R = 158.579760862474;   S = 0.006734972683744581; T = 0.006734972683744581 }
    @{ Row = 4;  E = 3; G = 91.47356466666666;  H = 274.420694;          I = 0.7914295280104694; J = 0.7914295280104694; K = 3; M = 0.572208;            N = 1.716624;           O = 0.02527946443038672;  P = 0.02527946443038671;  Q = 52.341905490784;     R = 471.077149417056;   S = 0.02000691460249841;  T = 0.02000691460249841  }
    @{ Row = 5;  E = 3; G = 20.94207233333333;  H = 62.82621699999999;   I = 0.181190866265331;  J = 0.181190866265331;  K = 3; M = 21.870458;           N = 65.61137400000001;  O = 0.9662106525726075;  P = 0.9662106525726075;  Q = 458.0127133991287;   R = 4122.114420592158;  S = 0.1750685451344215;  T = 0.1750685451344215  }
    @{ Row = 6;  E = 3; G = 20.94207233333333;  H = 62.82621699999999;   I = 0.181190866265331;  J = 0.181190866265331;  K = 3; M = 0.1926236666666667;  N = 0.577871;           O = 0.008509882997005752; P = 0.008509882997005752; Q = 4.033938760445222;   R = 36.305448844007;    S = 0.001541913072044084; T = 0.001541913072044084 }
    @{ Row = 7;  E = 3; G = 20.94207233333333;  H = 62.82621699999999;   I = 0.181190866265331;  J = 0.181190866265331;  K = 3; M = 0.572208;            N = 1.716624;           O = 0.02527946443038672;  P = 0.02527946443038671;  Q = 11.983221325712;     R = 107.848991931408;   S = 0.004580408058865392; T = 0.004580408058865391 }
    @{ Row = 8;  E = 3; G = 3.164539666666667;  H = 9.493619000000001;   I = 0.02737960572419959; J = 0.02737960572419959; K = 3; M = 21.870458;         N = 65.61137400000001;  O = 0.9662106525726075;  P = 0.9662106525726075;  Q = 69.20993186916735;   R = 622.8893868225061;  S = 0.02645446671395959;  T = 0.02645446671395959  }
    @{ Row = 9;  E = 3; G = 3.164539666666667;  H = 9.493619000000001;   I = 0.02737960572419959; J = 0.02737960572419959; K = 3; M = 0.1926236666666667; N = 0.577871;           O = 0.008509882997005752; P = 0.008509882997005752; Q = 0.6095652339054445;  R = 5.486087105149001;  S = 0.0002329972412170875; T = 0.0002329972412170875 }
    @{ Row = 10; E = 3; G = 3.164539666666667;  H = 9.493619000000001;   I = 0.02737960572419959; J = 0.02737960572419959; K = 3; M = 0.572208;          N = 1.716624;           O = 0.02527946443038672;  P = 0.02527946443038671;  Q = 1.810774913584;      R = 16.296974222256;    S = 0.0006921417690229162; T = 0.0006921417690229161 }
)

$cols = @("E", "G", "H", "I", "J", "K", "M", "N", "O", "P", "Q", "R", "S", "T")

foreach ($r in $rows) {
    foreach ($col in $cols) {
        $ws.Range("$col$($r.Row)").Value = $r[$col]
    }
}
